$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare formatting for the 3 brand-new rows (186-188) by copying the
#     format pattern (bold+border on col A, date format on col E) from row 185.
$ws.Range("A185:AC185").Copy()
$ws.Range("A186:AC186").PasteSpecial(-4122)
$ws.Range("A185:AC185").Copy()
$ws.Range("A187:AC187").PasteSpecial(-4122)
$ws.Range("A185:AC185").Copy()
$ws.Range("A188:AC188").PasteSpecial(-4122)

# row 185
$ws.Range("A185").Value = 183
$ws.Range("B185").Value = 7977862
$ws.Range("C185").Value = "Venezuela Primera Division"
$ws.Range("D185").Value = "Venezuela Primera Division"
$ws.Range("E185").Value = 45381.78125
$ws.Range("F185").Value = "Caracas"
$ws.Range("G185").Value = "Inter de Barinas"
$ws.Range("H185").Value = 1
$ws.Range("I185").Value = 2
$ws.Range("J185").Value = "A"
$ws.Range("K185").Value = 1.571
$ws.Range("L185").Value = 3.6
$ws.Range("M185").Value = 5.5
$ws.Range("N185").Value = 1.444
$ws.Range("O185").Value = 3.8
$ws.Range("P185").Value = 7
$ws.Range("Q185").Value = -1.25
$ws.Range("R185").Value = 2
$ws.Range("S185").Value = 1.8
$ws.Range("T185").Value = 2.25
$ws.Range("U185").Value = 1.9
$ws.Range("V185").Value = 1.9
$ws.Range("W185").Value = -1
$ws.Range("X185").Value = -1
$ws.Range("Y185").Value = 6
$ws.Range("Z185").Value = -1
$ws.Range("AA185").Value = 0.8
$ws.Range("AB185").Value = 0.8999999999999999
$ws.Range("AC185").Value = -1

# row 186
$ws.Range("A186").Value = 184
$ws.Range("B186").Value = 7977863
$ws.Range("C186").Value = "Venezuela Primera Division"
$ws.Range("D186").Value = "Venezuela Primera Division"
$ws.Range("E186").Value = 45381.89583333334
$ws.Range("F186").Value = "Portuguesa"
$ws.Range("G186").Value = "Carabobo"
$ws.Range("H186").Value = 0
$ws.Range("I186").Value = 0
$ws.Range("J186").Value = "D"
$ws.Range("K186").Value = 2.75
$ws.Range("L186").Value = 2.875
$ws.Range("M186").Value = 2.55
$ws.Range("N186").Value = 3
$ws.Range("O186").Value = 2.7
$ws.Range("P186").Value = 2.55
$ws.Range("Q186").Value = 0
$ws.Range("R186").Value = 2.05
$ws.Range("S186").Value = 1.75
$ws.Range("T186").Value = 2
$ws.Range("U186").Value = 1.975
$ws.Range("V186").Value = 1.825
$ws.Range("W186").Value = -1
$ws.Range("X186").Value = 1.7
$ws.Range("Y186").Value = -1
$ws.Range("Z186").Value = 0
$ws.Range("AA186").Value = 0
$ws.Range("AB186").Value = -1
$ws.Range("AC186").Value = 0.825

# row 187
$ws.Range("A187").Value = 185
$ws.Range("B187").Value = 7977380
$ws.Range("C187").Value = "Venezuela Primera Division"
$ws.Range("D187").Value = "Venezuela Primera Division"
$ws.Range("E187").Value = 45382.70833333334
$ws.Range("F187").Value = "Estudiantes Merida"
$ws.Range("G187").Value = "Deportivo La Guaira"
$ws.Range("H187").Value = 1
$ws.Range("I187").Value = 2
$ws.Range("J187").Value = "A"
$ws.Range("K187").Value = 3.4
$ws.Range("L187").Value = 2.875
$ws.Range("M187").Value = 2.15
$ws.Range("N187").Value = 3.1
$ws.Range("O187").Value = 2.8
$ws.Range("P187").Value = 2.375
$ws.Range("Q187").Value = 0.25
$ws.Range("R187").Value = 1.75
$ws.Range("S187").Value = 2.05
$ws.Range("T187").Value = 2.25
$ws.Range("U187").Value = 1.9
$ws.Range("V187").Value = 1.9
$ws.Range("W187").Value = -1
$ws.Range("X187").Value = -1
$ws.Range("Y187").Value = 1.375
$ws.Range("Z187").Value = -1
$ws.Range("AA187").Value = 1.05
$ws.Range("AB187").Value = 0.8999999999999999
$ws.Range("AC187").Value = -1

# row 188
$ws.Range("A188").Value = 186
$ws.Range("B188").Value = 8027687
$ws.Range("C188").Value = "Venezuela Primera Division"
$ws.Range("D188").Value = "Venezuela Primera Division"
$ws.Range("E188").Value = 45382.83333333334
$ws.Range("F188").Value = "Academia Puerto Cabello"
$ws.Range("G188").Value = "Monagas"
$ws.Range("H188").Value = 2
$ws.Range("I188").Value = 2
$ws.Range("J188").Value = "D"
$ws.Range("K188").Value = 1.75
$ws.Range("L188").Value = 3.4
$ws.Range("M188").Value = 4.333
$ws.Range("N188").Value = 1.85
$ws.Range("O188").Value = 3
$ws.Range("P188").Value = 4.2
$ws.Range("Q188").Value = -0.5
$ws.Range("R188").Value = 1.95
$ws.Range("S188").Value = 1.85
$ws.Range("T188").Value = 2
$ws.Range("U188").Value = 1.8
$ws.Range("V188").Value = 2
$ws.Range("W188").Value = -1
$ws.Range("X188").Value = 2
$ws.Range("Y188").Value = -1
$ws.Range("Z188").Value = -1
$ws.Range("AA188").Value = 0.8500000000000001
$ws.Range("AB188").Value = 0.8
$ws.Range("AC188").Value = -1

$ws.Range("B93").Value = 6236611
$ws.Range("F93").Value = "Mineros"
$ws.Range("G93").Value = "Monagas"
$ws.Range("H93").Value = 2
$ws.Range("I93").Value = 1
$ws.Range("J93").Value = "H"
$ws.Range("K93").Value = 3.2
$ws.Range("L93").Value = 3.4
$ws.Range("M93").Value = 2
$ws.Range("N93").Value = 4.2
$ws.Range("O93").Value = 3.8
$ws.Range("P93").Value = 1.65
$ws.Range("Q93").Value = 0.75
$ws.Range("R93").Value = 1.95
$ws.Range("S93").Value = 1.85
$ws.Range("T93").Value = 2.5
$ws.Range("U93").Value = 1.825
$ws.Range("V93").Value = 1.975
$ws.Range("W93").Value = 3.2
$ws.Range("X93").Value = -1
$ws.Range("Y93").Value = -1
$ws.Range("Z93").Value = 0.95
$ws.Range("AA93").Value = -1
$ws.Range("AB93").Value = 0.825
$ws.Range("AC93").Value = -1
$ws.Range("B94").Value = 6236612
$ws.Range("F94").Value = "Zamora"
$ws.Range("G94").Value = "Carabobo"
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 2
$ws.Range("J94").Value = "A"
$ws.Range("K94").Value = 3.2
$ws.Range("L94").Value = 3.1
$ws.Range("M94").Value = 2.15
$ws.Range("N94").Value = 4.5
$ws.Range("O94").Value = 3.3
$ws.Range("P94").Value = 1.75
$ws.Range("Q94").Value = 0.5
$ws.Range("R94").Value = 2
$ws.Range("S94").Value = 1.8
$ws.Range("T94").Value = 2.25
$ws.Range("U94").Value = 1.925
$ws.Range("V94").Value = 1.875
$ws.Range("W94").Value = -1
$ws.Range("X94").Value = -1
$ws.Range("Y94").Value = 0.75
$ws.Range("Z94").Value = -1
$ws.Range("AA94").Value = 0.8
$ws.Range("AB94").Value = -0.5
$ws.Range("AC94").Value = 0.4375
$ws.Range("B95").Value = 6236254
$ws.Range("F95").Value = "Academia Puerto Cabello"
$ws.Range("G95").Value = "Estudiantes Merida"
$ws.Range("H95").Value = 1
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = "H"
$ws.Range("K95").Value = 1.727
$ws.Range("L95").Value = 3.4
$ws.Range("M95").Value = 4.333
$ws.Range("N95").Value = 1.666
$ws.Range("O95").Value = 3.4
$ws.Range("P95").Value = 4.75
$ws.Range("Q95").Value = -0.75
$ws.Range("R95").Value = 1.875
$ws.Range("S95").Value = 1.925
$ws.Range("T95").Value = 2.5
$ws.Range("U95").Value = 1.9
$ws.Range("V95").Value = 1.9
$ws.Range("W95").Value = 0.6659999999999999
$ws.Range("X95").Value = -1
$ws.Range("Y95").Value = -1
$ws.Range("Z95").Value = 0.4375
$ws.Range("AA95").Value = -0.5
$ws.Range("AB95").Value = -1
$ws.Range("AC95").Value = 0.8999999999999999
$ws.Range("B96").Value = 6236255
$ws.Range("F96").Value = "Deportivo Rayo Zuliano"
$ws.Range("G96").Value = "Caracas"
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = "D"
$ws.Range("K96").Value = 3.75
$ws.Range("L96").Value = 3.1
$ws.Range("M96").Value = 1.95
$ws.Range("N96").Value = 2.9
$ws.Range("O96").Value = 2.875
$ws.Range("P96").Value = 2.45
$ws.Range("Q96").Value = 0.25
$ws.Range("R96").Value = 1.775
$ws.Range("S96").Value = 2.025
$ws.Range("T96").Value = 2.25
$ws.Range("U96").Value = 1.85
$ws.Range("V96").Value = 1.95
$ws.Range("W96").Value = -1
$ws.Range("X96").Value = 1.875
$ws.Range("Y96").Value = -1
$ws.Range("Z96").Value = 0.3875
$ws.Range("AA96").Value = -0.5
$ws.Range("AB96").Value = -1
$ws.Range("AC96").Value = 0.95
$ws.Range("B97").Value = 6236252
$ws.Range("F97").Value = "Deportivo Tachira"
$ws.Range("G97").Value = "CD Hermanos Colmenares"
$ws.Range("H97").Value = 1
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = "H"
$ws.Range("K97").Value = 1.363
$ws.Range("L97").Value = 4.2
$ws.Range("M97").Value = 7.5
$ws.Range("N97").Value = 1.333
$ws.Range("O97").Value = 4.5
$ws.Range("P97").Value = 8
$ws.Range("Q97").Value = -1.5
$ws.Range("R97").Value = 2
$ws.Range("S97").Value = 1.8
$ws.Range("T97").Value = 2.5
$ws.Range("U97").Value = 1.925
$ws.Range("V97").Value = 1.875
$ws.Range("W97").Value = 0.333
$ws.Range("X97").Value = -1
$ws.Range("Y97").Value = -1
$ws.Range("Z97").Value = -1
$ws.Range("AA97").Value = 0.8
$ws.Range("AB97").Value = -1
$ws.Range("AC97").Value = 0.875
$ws.Range("B98").Value = 6236251
$ws.Range("F98").Value = "Angostura FC"
$ws.Range("G98").Value = "Portuguesa"
$ws.Range("H98").Value = 1
$ws.Range("I98").Value = 2
$ws.Range("J98").Value = "A"
$ws.Range("K98").Value = 3.1
$ws.Range("L98").Value = 3.2
$ws.Range("M98").Value = 2.15
$ws.Range("N98").Value = 4
$ws.Range("O98").Value = 3.6
$ws.Range("P98").Value = 1.75
$ws.Range("Q98").Value = 0.75
$ws.Range("R98").Value = 1.8
$ws.Range("S98").Value = 2
$ws.Range("T98").Value = 2.5
$ws.Range("U98").Value = 1.95
$ws.Range("V98").Value = 1.85
$ws.Range("W98").Value = -1
$ws.Range("X98").Value = -1
$ws.Range("Y98").Value = 0.75
$ws.Range("Z98").Value = -0.5
$ws.Range("AA98").Value = 0.5
$ws.Range("AB98").Value = 0.95
$ws.Range("AC98").Value = -1

$ws.Range("B114").Value = 7352251
$ws.Range("F114").Value = "Caracas"
$ws.Range("G114").Value = "Academia Puerto Cabello"
$ws.Range("H114").Value = 1
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = "H"
$ws.Range("K114").Value = 2.1
$ws.Range("L114").Value = 3.2
$ws.Range("M114").Value = 3.3
$ws.Range("N114").Value = 2.15
$ws.Range("O114").Value = 3.1
$ws.Range("P114").Value = 3.2
$ws.Range("Q114").Value = -0.5
$ws.Range("R114").Value = 2.025
$ws.Range("S114").Value = 1.775
$ws.Range("T114").Value = 2.25
$ws.Range("U114").Value = 1.975
$ws.Range("V114").Value = 1.825
$ws.Range("W114").Value = 1.15
$ws.Range("X114").Value = -1
$ws.Range("Y114").Value = -1
$ws.Range("Z114").Value = 1.025
$ws.Range("AA114").Value = -1
$ws.Range("AB114").Value = -1
$ws.Range("AC114").Value = 0.825
$ws.Range("B115").Value = 7352250
$ws.Range("F115").Value = "Portuguesa"
$ws.Range("G115").Value = "Deportivo Tachira"
$ws.Range("H115").Value = 1
$ws.Range("I115").Value = 1
$ws.Range("J115").Value = "D"
$ws.Range("K115").Value = 3.1
$ws.Range("L115").Value = 2.875
$ws.Range("M115").Value = 2.3
$ws.Range("N115").Value = 3
$ws.Range("O115").Value = 2.875
$ws.Range("P115").Value = 2.375
$ws.Range("Q115").Value = 0.25
$ws.Range("R115").Value = 1.725
$ws.Range("S115").Value = 2.075
$ws.Range("T115").Value = 2
$ws.Range("U115").Value = 1.825
$ws.Range("V115").Value = 1.975
$ws.Range("W115").Value = -1
$ws.Range("X115").Value = 1.875
$ws.Range("Y115").Value = -1
$ws.Range("Z115").Value = 0.3625
$ws.Range("AA115").Value = -0.5
$ws.Range("AB115").Value = 0
$ws.Range("AC115").Value = 0
$ws.Range("B116").Value = 7352252
$ws.Range("F116").Value = "Deportivo Tachira"
$ws.Range("G116").Value = "Caracas"
$ws.Range("H116").Value = 1
$ws.Range("I116").Value = 1
$ws.Range("J116").Value = "D"
$ws.Range("K116").Value = 2.3
$ws.Range("L116").Value = 2.875
$ws.Range("M116").Value = 3.1
$ws.Range("N116").Value = 2.25
$ws.Range("O116").Value = 2.8
$ws.Range("P116").Value = 3.25
$ws.Range("Q116").Value = -0.25
$ws.Range("R116").Value = 1.975
$ws.Range("S116").Value = 1.825
$ws.Range("T116").Value = 2
$ws.Range("U116").Value = 1.925
$ws.Range("V116").Value = 1.875
$ws.Range("W116").Value = -1
$ws.Range("X116").Value = 1.8
$ws.Range("Y116").Value = -1
$ws.Range("Z116").Value = -0.5
$ws.Range("AA116").Value = 0.4125
$ws.Range("AB116").Value = 0
$ws.Range("AC116").Value = 0
$ws.Range("B117").Value = 7352254
$ws.Range("F117").Value = "Academia Puerto Cabello"
$ws.Range("G117").Value = "Portuguesa"
$ws.Range("H117").Value = 1
$ws.Range("I117").Value = 1
$ws.Range("J117").Value = "D"
$ws.Range("K117").Value = 2.05
$ws.Range("L117").Value = 3.4
$ws.Range("M117").Value = 3
$ws.Range("N117").Value = 1.833
$ws.Range("O117").Value = 3.5
$ws.Range("P117").Value = 3.5
$ws.Range("Q117").Value = -0.25
$ws.Range("R117").Value = 1.65
$ws.Range("S117").Value = 2.2
$ws.Range("T117").Value = 2.25
$ws.Range("U117").Value = 1.825
$ws.Range("V117").Value = 1.975
$ws.Range("W117").Value = -1
$ws.Range("X117").Value = 2.5
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = -0.5
$ws.Range("AA117").Value = 0.6000000000000001
$ws.Range("AB117").Value = -0.5
$ws.Range("AC117").Value = 0.4875

$ws.Range("B157").Value = 7920997
$ws.Range("F157").Value = "Carabobo"
$ws.Range("G157").Value = "UCV"
$ws.Range("H157").Value = 0
$ws.Range("I157").Value = 1
$ws.Range("J157").Value = "A"
$ws.Range("K157").Value = 1.833
$ws.Range("L157").Value = 3.1
$ws.Range("M157").Value = 4.2
$ws.Range("N157").Value = 1.833
$ws.Range("O157").Value = 3.1
$ws.Range("P157").Value = 4.2
$ws.Range("Q157").Value = -0.5
$ws.Range("R157").Value = 1.9
$ws.Range("S157").Value = 1.9
$ws.Range("T157").Value = 2
$ws.Range("U157").Value = 1.85
$ws.Range("V157").Value = 1.95
$ws.Range("W157").Value = -1
$ws.Range("X157").Value = -1
$ws.Range("Y157").Value = 3.2
$ws.Range("Z157").Value = -1
$ws.Range("AA157").Value = 0.8999999999999999
$ws.Range("AB157").Value = -1
$ws.Range("AC157").Value = 0.95
$ws.Range("B158").Value = 7920998
$ws.Range("F158").Value = "Zamora"
$ws.Range("G158").Value = "Caracas"
$ws.Range("H158").Value = 2
$ws.Range("I158").Value = 2
$ws.Range("J158").Value = "D"
$ws.Range("K158").Value = 3.75
$ws.Range("L158").Value = 3.2
$ws.Range("M158").Value = 1.909
$ws.Range("N158").Value = 3
$ws.Range("O158").Value = 2.9
$ws.Range("P158").Value = 2.375
$ws.Range("Q158").Value = 0.25
$ws.Range("R158").Value = 1.8
$ws.Range("S158").Value = 2
$ws.Range("T158").Value = 2
$ws.Range("U158").Value = 1.825
$ws.Range("V158").Value = 1.975
$ws.Range("W158").Value = -1
$ws.Range("X158").Value = 1.9
$ws.Range("Y158").Value = -1
$ws.Range("Z158").Value = 0.4
$ws.Range("AA158").Value = -0.5
$ws.Range("AB158").Value = 0.825
$ws.Range("AC158").Value = -1

